$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# C1 header text changes from "Text" to "WrappedText"
$ws.Range("C1").Value = "WrappedText"
# New D1 header "MultilineText" (same look as B1/C1: vertical-top, no wrap)
$ws.Range("D1").Value = "MultilineText"
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").WrapText = $false

# --- Update data row (row 2) ---
# Write D2 before C2 so the new shared-string table ends up in the same
# insertion order as the target workbook (MultilineText placeholder first).
# New D2 holds the multiline (non-wrapped) placeholder
$ws.Range("D2").Value = "{{item.MultilineText}}"
$ws.Range("D2").VerticalAlignment = -4160
$ws.Range("D2").WrapText = $false

# C2 keeps the wrapped-text style but now references the WrappedText placeholder
$ws.Range("C2").Value = "{{item.WrappedText}}"

# Set column D width to match the target layout (closest value reachable
# through this runtime's ColumnWidth rounding is used to land on the
# intended ~18.29-character-wide column)
$ws.Columns.Item(4).ColumnWidth = 17.5

# --- Update defined names to extend over the new column D ---
$wb.Names.Item("DataItems").RefersTo = "=Sheet1!`$A`$2:`$D`$3"
$wb.Names.Item("ReportHeaders").RefersTo = "=Sheet1!`$B`$1:`$D`$1"

# --- Update selection to match target (active cell D2, single cell selection) ---
$ws.Range("D2").Select()

$wb.Save()
